$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Items")

# New shared-string-backed cell values for row 10 (Buckler) and row 11 (Small Shield)
$ws.Range("K10").Value = "Ancient"
$ws.Range("L10").Value = 2
$ws.Range("M10").Value = 5
$ws.Range("N10").Value = "very high AR, low HP"
$ws.Range("N11").Value = "^ like RS barrows"

# New column widths for N (14) and O (15)
$ws.Columns.Item(14).ColumnWidth = 18.65
$ws.Columns.Item(15).ColumnWidth = 8.25

# Update view: move selection to N14 (also resets the scrolled topLeftCell back to default)
[void]$ws.Range("N14").Select()
